$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 111.583336
$ws.Range("I6").Value = 111.583336
$ws.Range("K6").Value = 334.750008
$ws.Range("M6").Value = -222.750008
$ws.Range("H28").Value = 244
$ws.Range("J28").Value = 2006
$ws.Range("L28").Value = 2006
$ws.Range("N28").Value = -2976
$ws.Range("H100").Value = 2356.4167
$ws.Range("I100").Value = 1740.7273
$ws.Range("K100").Value = 1740.7273
$ws.Range("M100").Value = -1199.7273
$ws.Range("H125").Value = 428.44446
$ws.Range("J125").Value = 444.8
$ws.Range("L125").Value = 4003.2
$ws.Range("N125").Value = -8923.200000000001
$ws.Range("H129").Value = 1569.3334
$ws.Range("I129").Value = 201.75
$ws.Range("J129").Value = 1757.9656
$ws.Range("K129").Value = 605.25
$ws.Range("L129").Value = 5273.8968
$ws.Range("M129").Value = 4394.75
$ws.Range("N129").Value = -15273.8968
$ws.Range("H132").Value = 55560576
$ws.Range("I132").Value = 71434700
$ws.Range("K132").Value = 214304100
$ws.Range("M132").Value = -214301570
$ws.Range("H137").Value = 92552.62
$ws.Range("I137").Value = 112559.586
$ws.Range("J137").Value = 2521.25
$ws.Range("K137").Value = 337678.758
$ws.Range("L137").Value = 7563.75
$ws.Range("M137").Value = -335128.758
$ws.Range("N137").Value = -12663.75
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1166.8
$ws.Range("I2").Value = 1108.5
$ws.Range("J2").Value = 1400
$ws.Range("K2").Value = 1108.5
$ws.Range("L2").Value = 1400
$ws.Range("M2").Value = -995.5
$ws.Range("N2").Value = -1626
$ws.Range("H32").Value = 4797.01
$ws.Range("I32").Value = 4113.8403
$ws.Range("J32").Value = 15500
$ws.Range("K32").Value = 4113.8403
$ws.Range("L32").Value = 15500
$ws.Range("M32").Value = -3826.8403
$ws.Range("N32").Value = -16074
$ws.Range("H61").Value = 7938565.5
$ws.Range("I61").Value = 10418147
$ws.Range("J61").Value = 3905.3
$ws.Range("K61").Value = 10418147
$ws.Range("L61").Value = 3905.3
$ws.Range("M61").Value = -10417935
$ws.Range("N61").Value = -4329.3
$ws.Range("H63").Value = 3908870
$ws.Range("I63").Value = 2993.3333
$ws.Range("J63").Value = 15626500
$ws.Range("K63").Value = 2993.3333
$ws.Range("L63").Value = 15626500
$ws.Range("M63").Value = -2307.3333
$ws.Range("N63").Value = -15627872
$ws.Range("H66").Value = 3908870
$ws.Range("I66").Value = 2993.3333
$ws.Range("J66").Value = 15626500
$ws.Range("K66").Value = 14966.6665
$ws.Range("L66").Value = 78132500
$ws.Range("M66").Value = -11534.6665
$ws.Range("N66").Value = -78139364
$ws.Range("H102").Value = 3505
$ws.Range("I102").Value = 3505
$ws.Range("J102").Value = 0
$ws.Range("K102").Value = 3505
$ws.Range("L102").Value = 0
$ws.Range("M102").Value = -1883
$ws.Range("N102").ClearContents()
$ws.Range("H116").Value = 1166.8
$ws.Range("I116").Value = 1108.5
$ws.Range("J116").Value = 1400
$ws.Range("K116").Value = 1108.5
$ws.Range("L116").Value = 1400
$ws.Range("M116").Value = 1185.5
$ws.Range("N116").Value = -5988
$ws.Range("H132").Value = 9445434
$ws.Range("I132").Value = 11112913
$ws.Range("J132").Value = 65869.25
$ws.Range("K132").Value = 33338739
$ws.Range("L132").Value = 197607.75
$ws.Range("M132").Value = -33336209
$ws.Range("N132").Value = -202667.75
$ws.Range("H136").Value = 7938565.5
$ws.Range("I136").Value = 10418147
$ws.Range("J136").Value = 3905.3
$ws.Range("K136").Value = 31254441
$ws.Range("L136").Value = 11715.9
$ws.Range("M136").Value = -31251891
$ws.Range("N136").Value = -16815.9
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1166.8
$ws.Range("I3").Value = 1108.5
$ws.Range("J3").Value = 1400
$ws.Range("K3").Value = 1108.5
$ws.Range("L3").Value = 1400
$ws.Range("M3").Value = -994.5
$ws.Range("N3").Value = -1628
$ws.Range("H86").Value = 2326.4517
$ws.Range("I86").Value = 2250.3
$ws.Range("J86").Value = 2464.9092
$ws.Range("K86").Value = 2250.3
$ws.Range("L86").Value = 2464.9092
$ws.Range("M86").Value = -1127.3
$ws.Range("N86").Value = -4710.9092
$ws.Range("H89").Value = 2326.4517
$ws.Range("I89").Value = 2250.3
$ws.Range("J89").Value = 2464.9092
$ws.Range("K89").Value = 11251.5
$ws.Range("L89").Value = 12324.546
$ws.Range("M89").Value = -5635.5
$ws.Range("N89").Value = -23556.546
$ws.Range("H108").Value = 43685
$ws.Range("J108").Value = 43685
$ws.Range("L108").Value = 43685
$ws.Range("N108").Value = -51365
$ws.Range("H134").Value = 3671.5715
$ws.Range("I134").Value = 3655.15
$ws.Range("K134").Value = 10965.45
$ws.Range("M134").Value = -8430.450000000001
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H55").Value = 13000
$ws.Range("J55").Value = 13000
$ws.Range("L55").Value = 13000
$ws.Range("N55").Value = -13630
$ws.Range("H107").Value = 2022.3846
$ws.Range("I107").Value = 740
$ws.Range("J107").Value = 2823.875
$ws.Range("K107").Value = 740
$ws.Range("L107").Value = 2823.875
$ws.Range("M107").Value = 1180
$ws.Range("N107").Value = -6663.875
$ws.Range("H109").Value = 165022830
$ws.Range("J109").Value = 165022830
$ws.Range("L109").Value = 165022830
$ws.Range("N109").Value = -165024910
$ws.Range("H132").Value = 62503684
$ws.Range("I132").Value = 76925520
$ws.Range("K132").Value = 230776560
$ws.Range("M132").Value = -230774030
$ws.Range("H134").Value = 125001090
$ws.Range("I134").Value = 125001090
$ws.Range("K134").Value = 375003270
$ws.Range("M134").Value = -375000735
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H87").Value = 16221.4
$ws.Range("I87").Value = 8119
$ws.Range("J87").Value = 28375
$ws.Range("K87").Value = 24357
$ws.Range("L87").Value = 85125
$ws.Range("M87").Value = -23109
$ws.Range("N87").Value = -87621
$ws.Range("H90").Value = 16221.4
$ws.Range("I90").Value = 8119
$ws.Range("J90").Value = 28375
$ws.Range("K90").Value = 73071
$ws.Range("L90").Value = 255375
$ws.Range("M90").Value = -66831
$ws.Range("N90").Value = -267855
$ws.Range("H109").Value = 3001.311
$ws.Range("I109").Value = 793.9167
$ws.Range("J109").Value = 3804
$ws.Range("K109").Value = 2381.7501
$ws.Range("L109").Value = 11412
$ws.Range("M109").Value = -1341.7501
$ws.Range("N109").Value = -13492
$ws.Range("H113").Value = 741.5454999999999
$ws.Range("I113").Value = 632.55554
$ws.Range("J113").Value = 817
$ws.Range("K113").Value = 1897.66662
$ws.Range("L113").Value = 2451
$ws.Range("M113").Value = 272.33338
$ws.Range("N113").Value = -6791
$ws.Range("H131").Value = 722.84
$ws.Range("J131").Value = 722.84
$ws.Range("L131").Value = 2168.52
$ws.Range("N131").Value = -12248.52
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H12").Value = 5908417
$ws.Range("I12").Value = 5908417
$ws.Range("K12").Value = 5908417
$ws.Range("M12").Value = -5908277
$ws.Range("H80").Value = 3793.25
$ws.Range("I80").Value = 3545.9092
$ws.Range("K80").Value = 3545.9092
$ws.Range("M80").Value = -2547.9092
$ws.Range("H83").Value = 3793.25
$ws.Range("I83").Value = 3545.9092
$ws.Range("K83").Value = 17729.546
$ws.Range("M83").Value = -12737.546
$ws.Range("H102").Value = 3527.9285
$ws.Range("I102").Value = 3442.0908
$ws.Range("J102").Value = 3842.6667
$ws.Range("K102").Value = 3442.0908
$ws.Range("L102").Value = 3842.6667
$ws.Range("M102").Value = -1820.0908
$ws.Range("N102").Value = -7086.6667
$ws.Range("H132").Value = 5104872.5
$ws.Range("I132").Value = 7942876
$ws.Range("J132").Value = 59533.11
$ws.Range("K132").Value = 23828628
$ws.Range("L132").Value = 178599.33
$ws.Range("M132").Value = -23826098
$ws.Range("N132").Value = -183659.33
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 2816
$ws.Range("I68").Value = 2450
$ws.Range("K68").Value = 2450
$ws.Range("M68").Value = -1701
$ws.Range("H71").Value = 2816
$ws.Range("I71").Value = 2450
$ws.Range("K71").Value = 12250
$ws.Range("M71").Value = -8506
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4420
$ws.Range("I62").Value = 3850
$ws.Range("J62").Value = 4800
$ws.Range("K62").Value = 3850
$ws.Range("L62").Value = 4800
$ws.Range("M62").Value = -3226
$ws.Range("N62").Value = -6048
$ws.Range("H65").Value = 4420
$ws.Range("I65").Value = 3850
$ws.Range("J65").Value = 4800
$ws.Range("K65").Value = 19250
$ws.Range("L65").Value = 24000
$ws.Range("M65").Value = -16130
$ws.Range("N65").Value = -30240
$ws.Range("H132").Value = 13158825
$ws.Range("I132").Value = 15625642
$ws.Range("J132").Value = 2467.3333
$ws.Range("K132").Value = 46876926
$ws.Range("L132").Value = 7401.999899999999
$ws.Range("M132").Value = -46874396
$ws.Range("N132").Value = -12461.9999
$ws.Range("H136").Value = 25179724
$ws.Range("I136").Value = 31281852
$ws.Range("J136").Value = 8450.625
$ws.Range("K136").Value = 93845556
$ws.Range("L136").Value = 25351.875
$ws.Range("M136").Value = -93843006
$ws.Range("N136").Value = -30451.875
